# Add new match-day rows to the "Partidos" sheet (Jun 21, 2025 / serial 45829)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Partidos")

$fecha = 45829

# Columns: A fecha | B jugador | C equipo | D posicion | E goles | F autogoles
#          G arquero(bool) | H goles_recibidos | I tarjetas_amarillas | J tarjetas_rojas
#          K asistencias | L Penales_Atajados
$rows = @(
    @("Gember Marin Sarria",        "Amarillo", "Arquero",       0, 1, $true,  4, 0, 0, 0, 0),
    @("Carlos Fernando Valencia",   "Amarillo", "Delantero",     1, 0, $false, 0, 0, 0, 0, 0),
    @("Sebastian Giraldo",          "Amarillo", "Mediocampista", 2, 0, $false, 0, 0, 0, 0, 0),
    @("Armando Murillo",            "Amarillo", "Defensa",       0, 0, $false, 0, 0, 0, 2, 0),
    @("Arnul David Narvaez",        "Amarillo", "Delantero",     0, 0, $false, 0, 1, 0, 0, 0),
    @("Fabian Caicedo",             "Azul",     "Arquero",       0, 0, $true,  3, 0, 0, 0, 0),
    @("Jairo Cuartas",              "Azul",     "Defensa",       1, 0, $false, 0, 0, 0, 0, 0),
    @("Edwin Hinestroza",           "Azul",     "Mediocampista", 1, 0, $false, 0, 0, 0, 0, 0),
    @("Julio Cesar Castaño",        "Azul",     "Mediocampista", 1, 0, $false, 0, 0, 0, 0, 0),
    @("Fabian Grajales",            "Azul",     "Mediocampista", 0, 0, $false, 0, 1, 0, 0, 0),
    @("Hermes Marquez",             "Azul",     "Defensa",       0, 0, $false, 0, 0, 0, 1, 0),
    @("Andres Tangarife",           "Azul",     "Delantero",     0, 0, $false, 0, 0, 0, 1, 0)
)

$startRow = 292
$lastExistingRow = $startRow - 1
$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($lastExistingRow, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $fecha
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
    $ws.Cells.Item($r, 9).Value = $row[7]
    $ws.Cells.Item($r, 10).Value = $row[8]
    $ws.Cells.Item($r, 11).Value = $row[9]
    $ws.Cells.Item($r, 12).Value = $row[10]
    $r = $r + 1
}

$lastRow = $r - 1
$excel.CutCopyMode = $false

# Re-establish the frozen header row and move the active selection to just
# past the newly appended data, like a user scrolling down after data entry.
$ws.Activate()
$win = $excel.Application.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$ws.Range("B305").Select()
